# Apply the "periods.xlsx" edit described in the commit:
# "add mask for columns, match columns name, rename columns in sav script"

$wb = $excel.ActiveWorkbook

$wsPeriod = $wb.Worksheets.Item("period_lbl")
$wsType   = $wb.Worksheets.Item("time_period_type")
$wsYear   = $wb.Worksheets.Item("year")

# --- period_lbl: append 4 new rows for the Jun-2024 period ---------------
$wsPeriod.Cells.Item(102, 1).Value = 101
$wsPeriod.Cells.Item(102, 2).Value = "2MATs: Jun 2024"
$wsPeriod.Cells.Item(104, 3).Value = "Month: 2024 (06) Jun"
$wsPeriod.Cells.Item(102, 3).Value = "2MATs: 2024 (06) Jun"
$wsPeriod.Cells.Item(103, 2).Value = "MAT: Jun 2024"
$wsPeriod.Cells.Item(103, 3).Value = "MAT: 2024 (06) Jun"
$wsPeriod.Cells.Item(104, 2).Value = "Month: Jun 2024"
$wsPeriod.Cells.Item(103, 1).Value = 102
$wsPeriod.Cells.Item(104, 1).Value = 103
$wsPeriod.Cells.Item(105, 2).Value = "3MMT: Jun 2024"
$wsPeriod.Cells.Item(105, 1).Value = 104
$wsPeriod.Cells.Item(105, 3).Value = "3MMT: 2024 (06) Jun"

# widen column B (period_lbl labels) so the longer text fits
$wsPeriod.Columns.Item(2).ColumnWidth = 21.1796875

# scroll period_lbl down to show the newly added rows, and leave the
# selection on the last new cell (B105)
[void]$wsPeriod.Select()
[void]$wsPeriod.Range("B105").Select()
$excel.ActiveWindow.ScrollRow = 98

# --- time_period_type: rename "2MAT/ 104 we" -> "2MATs/ 104 we" ----------
$wsType.Cells.Item(4, 2).Value = "2MATs/ 104 we"

# --- make "time_period_type" the active sheet/tab -------------------------
[void]$wsType.Select()
[void]$wsType.Range("C7").Select()

$wb.Save()
